# Insert a new "period_value" column right before the existing "aggregation"
# column (column AG), shifting "aggregation" one column to the right (to AH).
# The new column is populated with the constant 1 for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "aggregation" currently lives in column 33 (AG). Inserting a whole column
# there pushes it (and everything after it) one column to the right, to AH,
# exactly like choosing "Insert" on that column header in Excel.
$aggCol = 33

$ws.Columns.Item($aggCol).Insert()

# Header for the newly-inserted column.
$ws.Cells.Item(1, $aggCol).Value = "period_value"

# Fill every existing data row (2..157) with the constant value 1.
$lastRow = 157
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, $aggCol).Value = 1
}
